$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 45973
$ws.Range("D8").Value = 170.76
$ws.Range("E8").Value = 160.49
$ws.Range("F8").Value = 170.49
$ws.Range("G8").Value = 160.66

$ws.Range("A9").Value = 45973
$ws.Range("D9").Value = 170.76
$ws.Range("E9").Value = 160.49
$ws.Range("F9").Value = 170.49
$ws.Range("G9").Value = 160.66

$ws.Range("A10").Value = 45973
$ws.Range("D10").Value = 173.14
$ws.Range("E10").Value = 163.41999999999999
$ws.Range("F10").Value = 173.42
$ws.Range("G10").Value = 163.91

$ws.Range("A11").Value = 45972
$ws.Range("D11").Value = 169.58
$ws.Range("E11").Value = 160.28
$ws.Range("F11").Value = 170.28
$ws.Range("G11").Value = 160.44

$ws.Range("A12").Value = 45972
$ws.Range("D12").Value = 169.58
$ws.Range("E12").Value = 160.28
$ws.Range("F12").Value = 170.28
$ws.Range("G12").Value = 160.44

$ws.Range("A13").Value = 45972
$ws.Range("D13").Value = 172.63
$ws.Range("E13").Value = 163.22999999999999
$ws.Range("F13").Value = 173.23
$ws.Range("G13").Value = 163.72

$ws.Range("A17").Value = 45973
$ws.Range("D17").Value = 176.55
$ws.Range("E17").Value = 166.39
$ws.Range("F17").Value = 176.39

$ws.Range("A18").Value = 45972
$ws.Range("D18").Value = 176.07
$ws.Range("E18").Value = 166.2
$ws.Range("F18").Value = 176.2

$ws.Range("A22").Value = 45973
$ws.Range("D22").Value = 172
$ws.Range("E22").Value = 162.30000000000001
$ws.Range("F22").Value = 171.9
$ws.Range("G22").Value = 163.58000000000001

$ws.Range("A23").Value = 45973
$ws.Range("D23").Value = 177.92
$ws.Range("E23").Value = 167.14
$ws.Range("F23").Value = 177.14

$ws.Range("A24").Value = 45973
$ws.Range("D24").Value = 177.72
$ws.Range("E24").Value = 167.35
$ws.Range("F24").Value = 177.35

$ws.Range("A25").Value = 45973
$ws.Range("D25").Value = 178.55
$ws.Range("E25").Value = 166.75
$ws.Range("F25").Value = 176.75
$ws.Range("G25").Value = 166.79

$ws.Range("A26").Value = 45973
$ws.Range("D26").Value = 177.27
$ws.Range("E26").Value = 168.32
$ws.Range("F26").Value = 178.32

$ws.Range("A27").Value = 45972
$ws.Range("D27").Value = 171.59
$ws.Range("E27").Value = 162.22
$ws.Range("F27").Value = 171.82
$ws.Range("G27").Value = 163.51

$ws.Range("A28").Value = 45972
$ws.Range("D28").Value = 177.4
$ws.Range("E28").Value = 166.95
$ws.Range("F28").Value = 176.95

$ws.Range("A29").Value = 45972
$ws.Range("D29").Value = 177.21
$ws.Range("E29").Value = 167.15
$ws.Range("F29").Value = 177.15

$ws.Range("A30").Value = 45972
$ws.Range("D30").Value = 178.04
$ws.Range("E30").Value = 166.54
$ws.Range("F30").Value = 176.54
$ws.Range("G30").Value = 166.58

$ws.Range("A31").Value = 45972
$ws.Range("D31").Value = 176.76
$ws.Range("E31").Value = 168.11
$ws.Range("F31").Value = 178.11

$ws.Range("A35").Value = 45973
$ws.Range("D35").Value = 171.5
$ws.Range("E35").Value = 160.63
$ws.Range("F35").Value = 169.63

$ws.Range("A36").Value = 45972
$ws.Range("D36").Value = 170.98
$ws.Range("E36").Value = 160.44
$ws.Range("F36").Value = 169.44

$ws.Range("A40").Value = 45973
$ws.Range("D40").Value = 177.07
$ws.Range("E40").Value = 166.17
$ws.Range("F40").Value = 176.17

$ws.Range("A41").Value = 45973
$ws.Range("D41").Value = 176.77
$ws.Range("E41").Value = 166.59
$ws.Range("F41").Value = 176.59

$ws.Range("A42").Value = 45972
$ws.Range("D42").Value = 176.53
$ws.Range("E42").Value = 165.93
$ws.Range("F42").Value = 175.93

$ws.Range("A43").Value = 45972
$ws.Range("D43").Value = 176.24
$ws.Range("E43").Value = 166.35
$ws.Range("F43").Value = 176.35

$ws.Range("A47").Value = 45973
$ws.Range("D47").Value = 170.8
$ws.Range("E47").Value = 162.63999999999999
$ws.Range("F47").Value = 172.64

$ws.Range("A48").Value = 45973
$ws.Range("D48").Value = 170.8
$ws.Range("E48").Value = 162.81
$ws.Range("F48").Value = 172.81

$ws.Range("A49").Value = 45972
$ws.Range("D49").Value = 170.02
$ws.Range("E49").Value = 162.52000000000001
$ws.Range("F49").Value = 172.52

$ws.Range("A50").Value = 45972
$ws.Range("D50").Value = 170.02
$ws.Range("E50").Value = 162.69999999999999
$ws.Range("F50").Value = 172.7

$ws.Range("A54").Value = 45973
$ws.Range("D54").Value = 187.24
$ws.Range("E54").Value = 176.52
$ws.Range("F54").Value = 186.52

$ws.Range("A55").Value = 45973
$ws.Range("D55").Value = 174.89
$ws.Range("E55").Value = 173.76
$ws.Range("F55").Value = 183.76

$ws.Range("A56").Value = 45973
$ws.Range("D56").Value = 177.28

$ws.Range("A57").Value = 45973
$ws.Range("D57").Value = 176.88
$ws.Range("E57").Value = 168.03

$ws.Range("A58").Value = 45973
$ws.Range("D58").Value = 172.79
$ws.Range("E58").Value = 164.08
$ws.Range("F58").Value = 174.08

$ws.Range("A59").Value = 45973
$ws.Range("D59").Value = 179.53
$ws.Range("E59").Value = 174.66

$ws.Range("A60").Value = 45972
$ws.Range("D60").Value = 186.73
$ws.Range("E60").Value = 176.2
$ws.Range("F60").Value = 186.2

$ws.Range("A61").Value = 45972
$ws.Range("D61").Value = 174.38
$ws.Range("E61").Value = 173.71
$ws.Range("F61").Value = 183.71

$ws.Range("A62").Value = 45972
$ws.Range("D62").Value = 176.76

$ws.Range("A63").Value = 45972
$ws.Range("D63").Value = 176.41
$ws.Range("E63").Value = 167.98

$ws.Range("A64").Value = 45972
$ws.Range("D64").Value = 172.32
$ws.Range("E64").Value = 164.03
$ws.Range("F64").Value = 174.03

$ws.Range("A65").Value = 45972
$ws.Range("D65").Value = 179.05
$ws.Range("E65").Value = 174.4
